$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Normalize Spanish connector words to Title Case in place names ---
$ws.Range('B7').Value = 'Pabellón De Arteaga'
$ws.Range('B8').Value = 'Rincón De Romos'
$ws.Range('B31').Value = 'Amatenango De La Frontera'
$ws.Range('B32').Value = 'Amatenango Del Valle'
$ws.Range('B36').Value = 'Benemérito De Las Américas'
$ws.Range('B44').Value = 'Chiapa De Corzo'
$ws.Range('B49').Value = 'Comitán De Domínguez'
$ws.Range('B71').Value = 'Mazapa De Madero'
$ws.Range('B76').Value = 'Ocozocoautla De Espinosa'
$ws.Range('B85').Value = 'Salto De Agua'
$ws.Range('B86').Value = 'San Cristóbal De Las Casas'
$ws.Range('B123').Value = 'Guadalupe Y Calvo'
$ws.Range('B124').Value = 'Hidalgo Del Parral'
$ws.Range('B131').Value = 'San Francisco Del Oro'
$ws.Range('B134').Value = 'Valle De Zaragoza'
$ws.Range('B151').Value = 'San Juan De Sabinas'
$ws.Range('A164').Value = 'Ciudad De México'
$ws.Range('B168').Value = 'Cuajimalpa De Morelos'
$ws.Range('B192').Value = 'Nombre De Dios'
$ws.Range('B196').Value = 'Pánuco De Coronado'
$ws.Range('B202').Value = 'San Juan De Guadalupe'
$ws.Range('B203').Value = 'San Juan Del Río'
$ws.Range('A213').Value = 'Estado De México'
$ws.Range('B213').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B216').Value = 'Almoloya De Alquisiras'
$ws.Range('B217').Value = 'Almoloya De Juárez'
$ws.Range('B223').Value = 'Atizapán De Zaragoza'
$ws.Range('B230').Value = 'Coacalco De Berriozábal'
$ws.Range('B235').Value = 'Ecatepec De Morelos'
$ws.Range('B241').Value = 'Ixtapan De La Sal'
$ws.Range('B242').Value = 'Ixtapan Del Oro'
$ws.Range('B256').Value = 'Naucalpan De Juárez'
$ws.Range('B264').Value = 'San Antonio La Isla'
$ws.Range('B265').Value = 'San Felipe Del Progreso'
$ws.Range('B266').Value = 'San Simón De Guerrero'
$ws.Range('B268').Value = 'Soyaniquilpan De Juárez'
$ws.Range('B277').Value = 'Tenango Del Valle'
$ws.Range('B289').Value = 'Tlalnepantla De Baz'
$ws.Range('B295').Value = 'Valle De Bravo'
$ws.Range('B296').Value = 'Villa De Allende'
$ws.Range('B297').Value = 'Villa Del Carbón'
$ws.Range('B308').Value = 'San Miguel De Allende'
$ws.Range('B309').Value = 'Apaseo El Alto'
$ws.Range('B310').Value = 'Apaseo El Grande'
$ws.Range('B318').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B322').Value = 'Jaral Del Progreso'
$ws.Range('B330').Value = 'Purísima Del Rincón'
$ws.Range('B334').Value = 'San Diego De La Unión'
$ws.Range('B336').Value = 'San Francisco Del Rincón'
$ws.Range('B338').Value = 'San Luis De La Paz'
$ws.Range('B340').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B341').Value = 'Silao De La Victoria'
$ws.Range('B346').Value = 'Valle De Santiago'
$ws.Range('B352').Value = 'Acapulco De Juárez'
$ws.Range('B355').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B356').Value = 'Alcozauca De Guerrero'
$ws.Range('B360').Value = 'Atenango Del Río'
$ws.Range('B362').Value = 'Atoyac De Álvarez'
$ws.Range('B363').Value = 'Ayutla De Los Libres'
$ws.Range('B366').Value = 'Buenavista De Cuéllar'
$ws.Range('B367').Value = 'Chilapa De Álvarez'
$ws.Range('B368').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B369').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B374').Value = 'Coyuca De Benítez'
$ws.Range('B375').Value = 'Coyuca De Catalán'
$ws.Range('B379').Value = 'Cuetzala Del Progreso'
$ws.Range('B380').Value = 'Cutzamala De Pinzón'
$ws.Range('B386').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B387').Value = 'Iguala De La Independencia'
$ws.Range('B389').Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range('B390').Value = 'Zihuatanejo De Azueta'
$ws.Range('B392').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B395').Value = 'Mártir De Cuilapan'
$ws.Range('B408').Value = 'Taxco De Alarcón'
$ws.Range('B410').Value = 'Técpan De Galeana'
$ws.Range('B412').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B414').Value = 'Tixtla De Guerrero'
$ws.Range('B417').Value = 'Tlalixtaquilla De Maldonado'
$ws.Range('B418').Value = 'Tlapa De Comonfort'
$ws.Range('B430').Value = 'Agua Blanca De Iturbide'
$ws.Range('B435').Value = 'Atotonilco De Tula'
$ws.Range('B436').Value = 'Atotonilco El Grande'
$ws.Range('B442').Value = 'Cuautepec De Hinojosa'
$ws.Range('B445').Value = 'Huasca De Ocampo'
$ws.Range('B448').Value = 'Huejutla De Reyes'
$ws.Range('B451').Value = 'Jacala De Ledezma'
$ws.Range('B457').Value = 'Mineral Del Chico'
$ws.Range('B458').Value = 'Mineral Del Monte'
$ws.Range('B459').Value = 'Mixquiahuala De Juárez'
$ws.Range('B461').Value = 'Nopala De Villagrán'
$ws.Range('B462').Value = 'Omitlán De Juárez'
$ws.Range('B463').Value = 'Pachuca De Soto'
$ws.Range('B466').Value = 'Progreso De Obregón'
$ws.Range('B472').Value = 'Santiago De Anaya'
$ws.Range('B473').Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Range('B477').Value = 'Tenango De Doria'
$ws.Range('B479').Value = 'Tepehuacán De Guerrero'
$ws.Range('B480').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B483').Value = 'Tezontepec De Aldama'
$ws.Range('B490').Value = 'Tula De Allende'
$ws.Range('B491').Value = 'Tulancingo De Bravo'
$ws.Range('B494').Value = 'Zacualtipán De Ángeles'
$ws.Range('B495').Value = 'Zapotlán De Juárez'
$ws.Range('B500').Value = 'Ahualulco De Mercado'
$ws.Range('B503').Value = 'Atemajac De Brizuela'
$ws.Range('B504').Value = 'Atotonilco El Alto'
$ws.Range('B505').Value = 'Autlán De Navarro'
$ws.Range('B510').Value = 'Cañadas De Obregón'
$ws.Range('B515').Value = 'Cuautitlán De García Barragán'
$ws.Range('B521').Value = 'Encarnación De Díaz'
$ws.Range('B525').Value = 'Huejuquilla El Alto'
$ws.Range('B526').Value = 'Ixtlahuacán Del Río'
$ws.Range('B529').Value = 'Jilotlán De Los Dolores'
$ws.Range('B535').Value = 'Lagos De Moreno'
$ws.Range('B539').Value = 'Ojuelos De Jalisco'
$ws.Range('B547').Value = 'San Martín De Bolaños'
$ws.Range('B549').Value = 'San Miguel El Alto'
$ws.Range('B550').Value = 'San Sebastián Del Oeste'
$ws.Range('B551').Value = 'Santa María De Los Ángeles'
$ws.Range('B552').Value = 'Santa María Del Oro'
$ws.Range('B555').Value = 'Talpa De Allende'
$ws.Range('B556').Value = 'Tamazula De Gordiano'
$ws.Range('B558').Value = 'Techaluta De Montenegro'
$ws.Range('B562').Value = 'Teocuitatlán De Corona'
$ws.Range('B563').Value = 'Tepatitlán De Morelos'
$ws.Range('B565').Value = 'Tizapán El Alto'
$ws.Range('B566').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B571').Value = 'Unión De San Antonio'
$ws.Range('B572').Value = 'Unión De Tula'
$ws.Range('B573').Value = 'Valle De Guadalupe'
$ws.Range('B578').Value = 'Yahualica De González Gallo'
$ws.Range('B579').Value = 'Zacoalco De Torres'
$ws.Range('B581').Value = 'Zapotitlán De Vadillo'
$ws.Range('B582').Value = 'Zapotlán El Grande'
$ws.Range('B666').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B690').Value = 'Coatlán Del Río'
$ws.Range('B697').Value = 'Jonacatepec De Leandro Valle'
$ws.Range('B701').Value = 'Puente De Ixtla'
$ws.Range('B706').Value = 'Tetela Del Volcán'
$ws.Range('B707').Value = 'Tlaltizapán De Zapata'
$ws.Range('B713').Value = 'Zacualpan De Amilpas'
$ws.Range('B717').Value = 'Bahía De Banderas'
$ws.Range('B720').Value = 'Ixtlán Del Río'
$ws.Range('B726').Value = 'Santa María Del Oro'
$ws.Range('B741').Value = 'Lampazos De Naranjo'
$ws.Range('B743').Value = 'Mier Y Noriega'
$ws.Range('B747').Value = 'San Nicolás De Los Garza'
$ws.Range('B752').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B760').Value = 'Chalcatongo De Hidalgo'
$ws.Range('B762').Value = 'Coicoyán De Las Flores'
$ws.Range('B765').Value = 'Constancia Del Rosario'
$ws.Range('B768').Value = 'Cuyamecalco Villa De Zaragoza'
$ws.Range('B769').Value = 'Guevea De Humboldt'
$ws.Range('B770').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B771').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B772').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B773').Value = 'Ixtlán De Juárez'
$ws.Range('B774').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B784').Value = 'Mártires De Tacubaya'
$ws.Range('B787').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B790').Value = 'Nejapa De Madero'
$ws.Range('B791').Value = 'Oaxaca De Juárez'
$ws.Range('B792').Value = 'Ocotlán De Morelos'
$ws.Range('B793').Value = 'Pinotepa De Don Luis'
$ws.Range('B795').Value = 'Putla Villa De Guerrero'
$ws.Range('B806').Value = 'San Antonio De La Cal'
$ws.Range('B817').Value = 'San Felipe Jalapa De Díaz'
$ws.Range('B832').Value = 'San José Del Progreso'
$ws.Range('B839').Value = 'San Juan Bautista Lo De Soto'
$ws.Range('B846').Value = 'San Juan Del Estado'
$ws.Range('B847').Value = 'San Juan Del Río'
$ws.Range('B884').Value = 'San Miguel Del Puerto'
$ws.Range('B885').Value = 'San Miguel Del Río'
$ws.Range('B886').Value = 'San Miguel El Grande'
$ws.Range('B898').Value = 'San Pablo Villa De Mitla'
$ws.Range('B901').Value = 'San Pedro El Alto'
$ws.Range('B916').Value = 'San Pedro Y San Pablo Tequixtepec'
$ws.Range('B927').Value = 'Santa Ana Del Valle'
$ws.Range('B936').Value = 'Santa Cruz Tacache De Mina'
$ws.Range('B940').Value = 'Santa Inés Del Monte'
$ws.Range('B941').Value = 'Santa Lucía Del Camino'
$ws.Range('B951').Value = 'Santa María Del Tule'
$ws.Range('B957').Value = 'Santa María Jalapa Del Marqués'
$ws.Range('B959').Value = 'Santa María La Asunción'
$ws.Range('B982').Value = 'Santiago Del Río'
$ws.Range('B1011').Value = 'Santo Domingo De Morelos'
$ws.Range('B1029').Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range('B1031').Value = 'Tataltepec De Valdés'
$ws.Range('B1032').Value = 'Teotitlán De Flores Magón'
$ws.Range('B1033').Value = 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'
$ws.Range('B1034').Value = 'Tlacolula De Matamoros'
$ws.Range('B1036').Value = 'Totontepec Villa De Morelos'
$ws.Range('B1039').Value = 'Villa De Chilapa De Díaz'
$ws.Range('B1040').Value = 'Villa De Etla'
$ws.Range('B1041').Value = 'Villa De Tamazulápam Del Progreso'
$ws.Range('B1042').Value = 'Villa De Tututepec'
$ws.Range('B1043').Value = 'Villa De Zaachila'
$ws.Range('B1046').Value = 'Villa Sola De Vega'
$ws.Range('B1047').Value = 'Zapotitlán Del Río'
$ws.Range('B1050').Value = 'Zimatlán De Álvarez'
$ws.Range('B1068').Value = 'Chalchicomula De Sesma'
$ws.Range('B1076').Value = 'Chila De La Sal'
$ws.Range('B1092').Value = 'Huehuetlán El Chico'
$ws.Range('B1094').Value = 'Ixcamilpa De Guerrero'
$ws.Range('B1097').Value = 'Izúcar De Matamoros'
$ws.Range('B1110').Value = 'Palmar De Bravo'
$ws.Range('B1127').Value = 'San Salvador El Seco'
$ws.Range('B1128').Value = 'San Salvador El Verde'
$ws.Range('B1134').Value = 'Tecali De Herrera'
$ws.Range('B1141').Value = 'Tepanco De López'
$ws.Range('B1142').Value = 'Tepatlaxco De Hidalgo'
$ws.Range('B1148').Value = 'Tepeyahualco De Cuauhtémoc'
$ws.Range('B1149').Value = 'Tetela De Ocampo'
$ws.Range('B1154').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B1165').Value = 'Tuzamapan De Galeana'
$ws.Range('B1181').Value = 'Amealco De Bonfil'
$ws.Range('B1183').Value = 'Cadereyta De Montes'
$ws.Range('B1189').Value = 'Jalpan De Serra'
$ws.Range('B1190').Value = 'Landa De Matamoros'
$ws.Range('B1193').Value = 'Pinal De Amoles'
$ws.Range('B1196').Value = 'San Juan Del Río'
$ws.Range('B1207').Value = 'Armadillo De Los Infante'
$ws.Range('B1208').Value = 'Axtla De Terrazas'
$ws.Range('B1214').Value = 'Ciudad Del Maíz'
$ws.Range('B1223').Value = 'Mexquitic De Carmona'
$ws.Range('B1229').Value = 'San Ciro De Acosta'
$ws.Range('B1235').Value = 'Santa María Del Río'
$ws.Range('B1237').Value = 'Soledad De Graciano Sánchez'
$ws.Range('B1244').Value = 'Tanquián De Escobedo'
$ws.Range('B1248').Value = 'Villa De Arista'
$ws.Range('B1249').Value = 'Villa De Arriaga'
$ws.Range('B1250').Value = 'Villa De Guadalupe'
$ws.Range('B1251').Value = 'Villa De Ramos'
$ws.Range('B1252').Value = 'Villa De Reyes'
$ws.Range('B1295').Value = 'Jalpa De Méndez'
$ws.Range('B1331').Value = 'Soto La Marina'
$ws.Range('B1340').Value = 'Apetatitlán De Antonio Carvajal'
$ws.Range('B1349').Value = 'Nanacamilpa De Mariano Arista'
$ws.Range('B1352').Value = 'Papalotla De Xicohténcatl'
$ws.Range('B1353').Value = 'San Pablo Del Monte'
$ws.Range('B1354').Value = 'Sanctórum De Lázaro Cárdenas'
$ws.Range('B1358').Value = 'Tepetitla De Lardizábal'
$ws.Range('B1376').Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range('B1380').Value = 'Amatlán De Los Reyes'
$ws.Range('B1387').Value = 'Boca Del Río'
$ws.Range('B1389').Value = 'Camarón De Tejeda'
$ws.Range('B1392').Value = 'Castillo De Teayo'
$ws.Range('B1399').Value = 'Chinampa De Gorostiza'
$ws.Range('B1409').Value = 'Cosamaloapan De Carpio'
$ws.Range('B1410').Value = 'Cosautlán De Carvajal'
$ws.Range('B1426').Value = 'Hueyapan De Ocampo'
$ws.Range('B1427').Value = 'Huiloapan De Cuauhtémoc'
$ws.Range('B1428').Value = 'Ignacio De La Llave'
$ws.Range('B1432').Value = 'Ixhuacán De Los Reyes'
$ws.Range('B1433').Value = 'Ixhuatlán De Madero'
$ws.Range('B1434').Value = 'Ixhuatlán Del Café'
$ws.Range('B1435').Value = 'Ixhuatlán Del Sureste'
$ws.Range('B1446').Value = 'Juchique De Ferrer'
$ws.Range('B1449').Value = 'Lerdo De Tejada'
$ws.Range('B1452').Value = 'Martínez De La Torre'
$ws.Range('B1455').Value = 'Medellín De Bravo'
$ws.Range('B1458').Value = 'Mixtla De Altamirano'
$ws.Range('B1460').Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range('B1469').Value = 'Ozuluama De Mascareñas'
$ws.Range('B1473').Value = 'Paso De Ovejas'
$ws.Range('B1474').Value = 'Paso Del Macho'
$ws.Range('B1478').Value = 'Poza Rica De Hidalgo'
$ws.Range('B1489').Value = 'Sayula De Alemán'
$ws.Range('B1492').Value = 'Soledad De Doblado'
$ws.Range('B1498').Value = 'Tatahuicapan De Juárez'
$ws.Range('B1530').Value = 'Vega De Alatorre'
$ws.Range('B1552').Value = 'Cañitas De Felipe Pescador'
$ws.Range('B1554').Value = 'Concepción Del Oro'
$ws.Range('B1564').Value = 'Jiménez Del Teul'
$ws.Range('B1571').Value = 'Noria De Ángeles'
$ws.Range('B1581').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B1584').Value = 'Villa De Cos'

# --- Fix floating point precision on percentage values ---
$ws.Range("D4").Value = 0.000957723354768094
$ws.Range("D391").Value = 0.000957723354768094
$ws.Range("D405").Value = 0.000957723354768094
$ws.Range("D442").Value = 0.000957723354768094
$ws.Range("D517").Value = 0.000957723354768094
$ws.Range("D930").Value = 0.000957723354768094
$ws.Range("D1290").Value = 0.000957723354768094
$ws.Range("D1478").Value = 0.000957723354768094
$ws.Range("D312").Value = 0.009303598303461488

# --- Remove trailing footer/metadata rows (1593:1597) ---
$ws.Rows("1593:1597").Delete()
